# Fruta / hortaliza, semanal
# A new weekly price record (row) is inserted at row 29, pushing the
# existing rows 29-36 down to rows 30-37 (dimension grows from R36 to R37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 29.
$ws.Rows.Item(29).Insert()

# Fill in the new row 29 with the new weekly record.
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(29, 3).Value = "Bíobío"
$ws.Cells.Item(29, 4).Value = 44636
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(29, 6).Value = 100112030
$ws.Cells.Item(29, 7).Value = "Poroto granado"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 80
$ws.Cells.Item(29, 11).Value = 22000
$ws.Cells.Item(29, 12).Value = 23000
$ws.Cells.Item(29, 13).Value = 22375
$ws.Cells.Item(29, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29, 15).Value = "Región Metropolitana"
$ws.Cells.Item(29, 16).Value = 895
$ws.Cells.Item(29, 17).Value = 25
$ws.Cells.Item(29, 18).Value = "Hortaliza"
